$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1288.75
$ws.Range("I6").Value = 1708.1666
$ws.Range("J6").Value = 30.5
$ws.Range("K6").Value = 5124.4998
$ws.Range("L6").Value = 91.5
$ws.Range("M6").Value = -5012.4998
$ws.Range("N6").Value = -315.5
$ws.Range("H8").Value = 57.333332
$ws.Range("I8").Value = 57.333332
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 171.999996
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -32.99999600000001
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("H39").Value = 117.6
$ws.Range("I39").Value = 117.6
$ws.Range("K39").Value = 352.8
$ws.Range("M39").Value = -56.79999999999995
$ws.Range("H58").Value = 401
$ws.Range("I58").Value = 401
$ws.Range("K58").Value = 1203
$ws.Range("M58").Value = -1053
$ws.Range("H137").Value = 1698.6909
$ws.Range("J137").Value = 2321.3333
$ws.Range("L137").Value = 6963.999899999999
$ws.Range("N137").Value = -12063.9999
$ws.Range("H138").Value = 4782.4688
$ws.Range("J138").Value = 5057.92
$ws.Range("L138").Value = 15173.76
$ws.Range("N138").Value = -25453.76
$ws.Range("N8").ClearContents()
$ws.Range("M31").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17787.912
$ws.Range("I32").Value = 17577.334
$ws.Range("K32").Value = 17577.334
$ws.Range("M32").Value = -17290.334
$ws.Range("H45").Value = 3997.6667
$ws.Range("I45").Value = 3996.75
$ws.Range("K45").Value = 3996.75
$ws.Range("M45").Value = -3619.75
$ws.Range("H61").Value = 1702.5
$ws.Range("I61").Value = 1506.6666
$ws.Range("K61").Value = 1506.6666
$ws.Range("M61").Value = -1294.6666
$ws.Range("H74").Value = 1702.8966
$ws.Range("I74").Value = 1141.8182
$ws.Range("J74").Value = 3466.2856
$ws.Range("K74").Value = 1141.8182
$ws.Range("L74").Value = 3466.2856
$ws.Range("M74").Value = -267.8181999999999
$ws.Range("N74").Value = -5214.2856
$ws.Range("H77").Value = 1702.8966
$ws.Range("I77").Value = 1141.8182
$ws.Range("J77").Value = 3466.2856
$ws.Range("K77").Value = 5709.090999999999
$ws.Range("L77").Value = 17331.428
$ws.Range("M77").Value = -1341.090999999999
$ws.Range("N77").Value = -26067.428
$ws.Range("H102").Value = 1371.75
$ws.Range("I102").Value = 1371.75
$ws.Range("K102").Value = 1371.75
$ws.Range("M102").Value = 250.25
$ws.Range("H132").Value = 2943.75
$ws.Range("J132").Value = 3628.6667
$ws.Range("L132").Value = 10886.0001
$ws.Range("N132").Value = -15946.0001
$ws.Range("H133").Value = 76685
$ws.Range("J133").Value = 76685
$ws.Range("L133").Value = 76685
$ws.Range("N133").Value = -81745
$ws.Range("H136").Value = 1702.5
$ws.Range("I136").Value = 1506.6666
$ws.Range("K136").Value = 4519.9998
$ws.Range("M136").Value = -1969.9998
$ws.Range("H138").Value = 99997
$ws.Range("J138").Value = 99997
$ws.Range("L138").Value = 99997
$ws.Range("N138").Value = -110277

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10648.5
$ws.Range("I20").Value = 9865
$ws.Range("K20").Value = 9865
$ws.Range("M20").Value = -9618
$ws.Range("H81").Value = 74973.336
$ws.Range("J81").Value = 74973.336
$ws.Range("L81").Value = 74973.336
$ws.Range("N81").Value = -77095.336
$ws.Range("H84").Value = 74973.336
$ws.Range("J84").Value = 74973.336
$ws.Range("L84").Value = 224920.008
$ws.Range("N84").Value = -235528.008
$ws.Range("H107").Value = 1068.85
$ws.Range("I107").Value = 904.7059
$ws.Range("K107").Value = 904.7059
$ws.Range("M107").Value = 1015.2941
$ws.Range("H134").Value = 3466.2856
$ws.Range("I134").Value = 3529.4614
$ws.Range("J134").Value = 3363.625
$ws.Range("K134").Value = 10588.3842
$ws.Range("L134").Value = 10090.875
$ws.Range("M134").Value = -8053.3842
$ws.Range("N134").Value = -15160.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2092.2
$ws.Range("I31").Value = 2054.4583
$ws.Range("J31").Value = 2998
$ws.Range("K31").Value = 2054.4583
$ws.Range("L31").Value = 2998
$ws.Range("M31").Value = -1759.4583
$ws.Range("N31").Value = -3588
$ws.Range("H34").Value = 2092.2
$ws.Range("I34").Value = 2054.4583
$ws.Range("J34").Value = 2998
$ws.Range("K34").Value = 2054.4583
$ws.Range("L34").Value = 2998
$ws.Range("M34").Value = -1852.4583
$ws.Range("N34").Value = -3402
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H58").Value = 3028.6667
$ws.Range("I58").Value = 3028.6667
$ws.Range("K58").Value = 3028.6667
$ws.Range("M58").Value = -2825.6667
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("H60").Value = 36376.953
$ws.Range("J60").Value = 36871.094
$ws.Range("L60").Value = 36871.094
$ws.Range("N60").Value = -37893.094
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H99").Value = 9391.714
$ws.Range("I99").Value = 9185.75
$ws.Range("J99").Value = 9666.333000000001
$ws.Range("K99").Value = 9185.75
$ws.Range("L99").Value = 9666.333000000001
$ws.Range("M99").Value = -7687.75
$ws.Range("N99").Value = -12662.333
$ws.Range("H126").Value = 9391.714
$ws.Range("I126").Value = 9185.75
$ws.Range("J126").Value = 9666.333000000001
$ws.Range("K126").Value = 27557.25
$ws.Range("L126").Value = 28998.999
$ws.Range("M126").Value = -25087.25
$ws.Range("N126").Value = -33938.999
$ws.Range("H136").Value = 3028.6667
$ws.Range("I136").Value = 3028.6667
$ws.Range("K136").Value = 9086.000100000001
$ws.Range("M136").Value = -6536.000100000001
$ws.Range("H141").Value = 60747.535
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 60747.535
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 60747.535
$ws.Range("N141").Value = -71107.535
$ws.Range("N50").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("N59").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("M141").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 361.66666
$ws.Range("I33").Value = 384.5
$ws.Range("J33").Value = 316
$ws.Range("K33").Value = 2307
$ws.Range("L33").Value = 1896
$ws.Range("M33").Value = -2024
$ws.Range("N33").Value = -2462
$ws.Range("H46").Value = 3466
$ws.Range("I46").Value = 3199
$ws.Range("K46").Value = 9597
$ws.Range("M46").Value = -9506
$ws.Range("H103").Value = 994.7143
$ws.Range("I103").Value = 762
$ws.Range("J103").Value = 1305
$ws.Range("K103").Value = 2286
$ws.Range("L103").Value = 3915
$ws.Range("M103").Value = -1407
$ws.Range("N103").Value = -5673
$ws.Range("H134").Value = 1533.8572
$ws.Range("I134").Value = 1533.8572
$ws.Range("K134").Value = 4601.571599999999
$ws.Range("M134").Value = 468.4284000000007

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2007
$ws.Range("I122").Value = 2007
$ws.Range("K122").Value = 6021
$ws.Range("M122").Value = -3571

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1600.75
$ws.Range("I7").Value = 1600.75
$ws.Range("K7").Value = 1600.75
$ws.Range("M7").Value = -1488.75
$ws.Range("H55").Value = 388.7143
$ws.Range("I55").Value = 384.4
$ws.Range("K55").Value = 384.4
$ws.Range("M55").Value = -211.4
$ws.Range("H68").Value = 1784
$ws.Range("J68").Value = 2099.3333
$ws.Range("L68").Value = 2099.3333
$ws.Range("N68").Value = -3597.3333
$ws.Range("H71").Value = 1784
$ws.Range("J71").Value = 2099.3333
$ws.Range("L71").Value = 10496.6665
$ws.Range("N71").Value = -17984.6665
$ws.Range("H100").Value = 2033.875
$ws.Range("I100").Value = 1305.5
$ws.Range("J100").Value = 2762.25
$ws.Range("K100").Value = 1305.5
$ws.Range("L100").Value = 2762.25
$ws.Range("M100").Value = -764.5
$ws.Range("N100").Value = -3844.25
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 1600.75
$ws.Range("I126").Value = 1600.75
$ws.Range("K126").Value = 4802.25
$ws.Range("M126").Value = -2332.25
$ws.Range("H132").Value = 5999
$ws.Range("J132").Value = 5999
$ws.Range("L132").Value = 17997
$ws.Range("N132").Value = -23057
$ws.Range("N125").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2997.9
$ws.Range("I122").Value = 3155.1428
$ws.Range("J122").Value = 2631
$ws.Range("K122").Value = 9465.428400000001
$ws.Range("L122").Value = 7893
$ws.Range("M122").Value = -7015.428400000001
$ws.Range("N122").Value = -12793
